$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in column H1, matching the formatting of the
# existing header cells (e.g. G1: bold, centered, bordered style)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Add the corresponding numeric value for the data row, unstyled like
# the neighboring numeric cells (F2/G2)
$ws.Range("H2").Value = 0
